# TradingModel - 2021/11/12 data update
# Align A4's date formatting with the other data rows (YYYY-MM-DD HH:MM:SS)
# and append the new day's data (2021-11-12) in row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 switches from the "date only" format to the "date + time" format
# used by the rows above it.
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row for 2021-11-12.
$ws.Range("A5").Value = 44512
$ws.Range("A5").NumberFormat = "YYYY-MM-DD"
$ws.Range("B5").Value = -1888.2
